$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows at 670-671 (pushes the existing 670:711 block down to 672:713)
$ws.Rows("670:671").Insert()

# Row 670: 2026/01/17 (Sat)
$ws.Range("A670").NumberFormat = "@"
$ws.Range("A670").Value = "2026/01/17"
$ws.Range("A670").Style = "Normal"
$ws.Range("B670").Value = "土"
$ws.Range("C670").Value = 22
$ws.Range("D670").Value = 201

# Row 671: 2026/01/18 (Sun)
$ws.Range("A671").NumberFormat = "@"
$ws.Range("A671").Value = "2026/01/18"
$ws.Range("A671").Style = "Normal"
$ws.Range("B671").Value = "日"
$ws.Range("C671").Value = 2
$ws.Range("D671").Value = 170
